$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 3010.5557
$ws.Cells.Item(19, 9).Value = 2749
$ws.Cells.Item(19, 11).Value = 2749
$ws.Cells.Item(19, 13).Value = -2574
$ws.Cells.Item(93, 8).Value = 49646.855
$ws.Cells.Item(93, 10).Value = 49646.855
$ws.Cells.Item(93, 12).Value = 49646.855
$ws.Cells.Item(93, 14).Value = -54638.855
$ws.Cells.Item(116, 8).Value = 31255964
$ws.Cells.Item(116, 9).Value = 83337830
$ws.Cells.Item(116, 10).Value = 6841.2
$ws.Cells.Item(116, 11).Value = 83337830
$ws.Cells.Item(116, 12).Value = 6841.2
$ws.Cells.Item(116, 13).Value = -83334388
$ws.Cells.Item(116, 14).Value = -13725.2
$ws.Cells.Item(132, 8).Value = 2053.476
$ws.Cells.Item(132, 9).Value = 2058.5264
$ws.Cells.Item(132, 11).Value = 6175.5792
$ws.Cells.Item(132, 13).Value = -3645.5792
$ws.Cells.Item(141, 8).Value = 12822609
$ws.Cells.Item(141, 9).Value = 13335459
$ws.Cells.Item(141, 11).Value = 40006377
$ws.Cells.Item(141, 13).Value = -40001197

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1292249.4
$ws.Cells.Item(32, 9).Value = 1376626.5
$ws.Cells.Item(32, 10).Value = 12529.167
$ws.Cells.Item(32, 11).Value = 1376626.5
$ws.Cells.Item(32, 12).Value = 12529.167
$ws.Cells.Item(32, 13).Value = -1376339.5
$ws.Cells.Item(32, 14).Value = -13103.167
$ws.Cells.Item(57, 8).Value = 4527
$ws.Cells.Item(57, 9).Value = 4527
$ws.Cells.Item(57, 11).Value = 4527
$ws.Cells.Item(57, 13).Value = -4043
$ws.Cells.Item(74, 8).Value = 37349.574
$ws.Cells.Item(74, 9).Value = 49680.324
$ws.Cells.Item(74, 11).Value = 49680.324
$ws.Cells.Item(74, 13).Value = -48806.324
$ws.Cells.Item(77, 8).Value = 37349.574
$ws.Cells.Item(77, 9).Value = 49680.324
$ws.Cells.Item(77, 11).Value = 248401.62
$ws.Cells.Item(77, 13).Value = -244033.62
$ws.Cells.Item(119, 8).Value = 73479
$ws.Cells.Item(119, 10).Value = 73479
$ws.Cells.Item(119, 12).Value = 73479
$ws.Cells.Item(119, 14).Value = -83155
$ws.Cells.Item(126, 8).Value = 5241.875
$ws.Cells.Item(126, 9).Value = 5241.875
$ws.Cells.Item(126, 11).Value = 15725.625
$ws.Cells.Item(126, 13).Value = -13255.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 37078680
$ws.Cells.Item(86, 9).Value = 62749.824
$ws.Cells.Item(86, 11).Value = 62749.824
$ws.Cells.Item(86, 13).Value = -61626.824
$ws.Cells.Item(89, 8).Value = 37078680
$ws.Cells.Item(89, 9).Value = 62749.824
$ws.Cells.Item(89, 11).Value = 313749.12
$ws.Cells.Item(89, 13).Value = -308133.12
$ws.Cells.Item(113, 8).Value = 5036.125
$ws.Cells.Item(113, 9).Value = 5036.125
$ws.Cells.Item(113, 11).Value = 5036.125
$ws.Cells.Item(113, 13).Value = -2866.125
$ws.Cells.Item(134, 8).Value = 4660.0186
$ws.Cells.Item(134, 9).Value = 1263.1666
$ws.Cells.Item(134, 10).Value = 8906.083000000001
$ws.Cells.Item(134, 11).Value = 3789.4998
$ws.Cells.Item(134, 12).Value = 26718.249
$ws.Cells.Item(134, 13).Value = -1254.4998
$ws.Cells.Item(134, 14).Value = -31788.249

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 6977.5435
$ws.Cells.Item(31, 9).Value = 2522.5908
$ws.Cells.Item(31, 10).Value = 11061.25
$ws.Cells.Item(31, 11).Value = 2522.5908
$ws.Cells.Item(31, 12).Value = 11061.25
$ws.Cells.Item(31, 13).Value = -2227.5908
$ws.Cells.Item(31, 14).Value = -11651.25
$ws.Cells.Item(34, 8).Value = 6977.5435
$ws.Cells.Item(34, 9).Value = 2522.5908
$ws.Cells.Item(34, 10).Value = 11061.25
$ws.Cells.Item(34, 11).Value = 2522.5908
$ws.Cells.Item(34, 12).Value = 11061.25
$ws.Cells.Item(34, 13).Value = -2320.5908
$ws.Cells.Item(34, 14).Value = -11465.25
$ws.Cells.Item(62, 8).Value = 15631371
$ws.Cells.Item(62, 9).Value = 31256492
$ws.Cells.Item(62, 11).Value = 31256492
$ws.Cells.Item(62, 13).Value = -31255868
$ws.Cells.Item(64, 8).Value = 67490.836
$ws.Cells.Item(64, 10).Value = 67490.836
$ws.Cells.Item(64, 12).Value = 67490.836
$ws.Cells.Item(64, 14).Value = -67986.836
$ws.Cells.Item(65, 8).Value = 15631371
$ws.Cells.Item(65, 9).Value = 31256492
$ws.Cells.Item(65, 11).Value = 156282460
$ws.Cells.Item(65, 13).Value = -156279340
$ws.Cells.Item(67, 8).Value = 67490.836
$ws.Cells.Item(67, 10).Value = 67490.836
$ws.Cells.Item(67, 12).Value = 67490.836
$ws.Cells.Item(67, 14).Value = -69206.836
$ws.Cells.Item(69, 8).Value = 31998.75
$ws.Cells.Item(69, 9).Value = 28331.666
$ws.Cells.Item(69, 10).Value = 43000
$ws.Cells.Item(69, 11).Value = 28331.666
$ws.Cells.Item(69, 12).Value = 43000
$ws.Cells.Item(69, 13).Value = -27582.666
$ws.Cells.Item(69, 14).Value = -44498
$ws.Cells.Item(72, 8).Value = 31998.75
$ws.Cells.Item(72, 9).Value = 28331.666
$ws.Cells.Item(72, 10).Value = 43000
$ws.Cells.Item(72, 11).Value = 84994.99800000001
$ws.Cells.Item(72, 12).Value = 129000
$ws.Cells.Item(72, 13).Value = -81250.99800000001
$ws.Cells.Item(72, 14).Value = -136488
$ws.Cells.Item(76, 8).Value = 4953.857
$ws.Cells.Item(76, 9).Value = 4953.857
$ws.Cells.Item(76, 11).Value = 4953.857
$ws.Cells.Item(76, 13).Value = -4638.857
$ws.Cells.Item(79, 8).Value = 4953.857
$ws.Cells.Item(79, 9).Value = 4953.857
$ws.Cells.Item(79, 11).Value = 4953.857
$ws.Cells.Item(79, 13).Value = -3861.857
$ws.Cells.Item(99, 8).Value = 12823.556
$ws.Cells.Item(99, 9).Value = 16853
$ws.Cells.Item(99, 10).Value = 9600
$ws.Cells.Item(99, 11).Value = 16853
$ws.Cells.Item(99, 12).Value = 9600
$ws.Cells.Item(99, 13).Value = -15355
$ws.Cells.Item(99, 14).Value = -12596
$ws.Cells.Item(126, 8).Value = 12823.556
$ws.Cells.Item(126, 9).Value = 16853
$ws.Cells.Item(126, 10).Value = 9600
$ws.Cells.Item(126, 11).Value = 50559
$ws.Cells.Item(126, 12).Value = 28800
$ws.Cells.Item(126, 13).Value = -48089
$ws.Cells.Item(126, 14).Value = -33740
$ws.Cells.Item(132, 8).Value = 5255.525
$ws.Cells.Item(132, 9).Value = 2193.5557
$ws.Cells.Item(132, 11).Value = 6580.6671
$ws.Cells.Item(132, 13).Value = -4050.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(11, 8).Value = 1279143.4
$ws.Cells.Item(11, 9).Value = 1643728.6
$ws.Cells.Item(11, 10).Value = 3095
$ws.Cells.Item(11, 11).Value = 4931185.800000001
$ws.Cells.Item(11, 12).Value = 9285
$ws.Cells.Item(11, 13).Value = -4931045.800000001
$ws.Cells.Item(11, 14).Value = -9565
$ws.Cells.Item(38, 8).Value = 124.71429
$ws.Cells.Item(38, 10).Value = 96.666664
$ws.Cells.Item(38, 12).Value = 289.999992
$ws.Cells.Item(38, 14).Value = -983.999992
$ws.Cells.Item(137, 8).Value = 113025.055
$ws.Cells.Item(137, 9).Value = 92486.27
$ws.Cells.Item(137, 10).Value = 145300.28
$ws.Cells.Item(137, 11).Value = 277458.81
$ws.Cells.Item(137, 12).Value = 435900.84
$ws.Cells.Item(137, 13).Value = -272358.81
$ws.Cells.Item(137, 14).Value = -446100.84

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 3635.8125
$ws.Cells.Item(102, 9).Value = 3958.6
$ws.Cells.Item(102, 10).Value = 3097.8333
$ws.Cells.Item(102, 11).Value = 3958.6
$ws.Cells.Item(102, 12).Value = 3097.8333
$ws.Cells.Item(102, 13).Value = -2336.6
$ws.Cells.Item(102, 14).Value = -6341.8333
$ws.Cells.Item(113, 8).Value = 6103.3887
$ws.Cells.Item(113, 9).Value = 3021.4666
$ws.Cells.Item(113, 10).Value = 8304.762000000001
$ws.Cells.Item(113, 11).Value = 3021.4666
$ws.Cells.Item(113, 12).Value = 8304.762000000001
$ws.Cells.Item(113, 13).Value = -851.4666000000002
$ws.Cells.Item(113, 14).Value = -12644.762

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1530
$ws.Cells.Item(22, 9).Value = 973.3333
$ws.Cells.Item(22, 10).Value = 2573.75
$ws.Cells.Item(22, 11).Value = 973.3333
$ws.Cells.Item(22, 12).Value = 2573.75
$ws.Cells.Item(22, 13).Value = -678.3333
$ws.Cells.Item(22, 14).Value = -3163.75
$ws.Cells.Item(27, 8).Value = 1530
$ws.Cells.Item(27, 9).Value = 973.3333
$ws.Cells.Item(27, 10).Value = 2573.75
$ws.Cells.Item(27, 11).Value = 973.3333
$ws.Cells.Item(27, 12).Value = 2573.75
$ws.Cells.Item(27, 13).Value = -866.3333
$ws.Cells.Item(27, 14).Value = -2787.75
$ws.Cells.Item(46, 8).Value = 1918045.8
$ws.Cells.Item(46, 9).Value = 6897070
$ws.Cells.Item(46, 10).Value = 3036.5386
$ws.Cells.Item(46, 11).Value = 6897070
$ws.Cells.Item(46, 12).Value = 3036.5386
$ws.Cells.Item(46, 13).Value = -6896882
$ws.Cells.Item(46, 14).Value = -3412.5386
$ws.Cells.Item(82, 8).Value = 2640.8572
$ws.Cells.Item(82, 9).Value = 805.5
$ws.Cells.Item(82, 10).Value = 3375
$ws.Cells.Item(82, 11).Value = 805.5
$ws.Cells.Item(82, 12).Value = 3375
$ws.Cells.Item(82, 13).Value = -444.5
$ws.Cells.Item(82, 14).Value = -4097
$ws.Cells.Item(85, 8).Value = 2640.8572
$ws.Cells.Item(85, 9).Value = 805.5
$ws.Cells.Item(85, 10).Value = 3375
$ws.Cells.Item(85, 11).Value = 805.5
$ws.Cells.Item(85, 12).Value = 3375
$ws.Cells.Item(85, 13).Value = 442.5
$ws.Cells.Item(85, 14).Value = -5871
$ws.Cells.Item(136, 8).Value = 8767.781000000001
$ws.Cells.Item(136, 9).Value = 2112.4285
$ws.Cells.Item(136, 11).Value = 6337.2855
$ws.Cells.Item(136, 13).Value = -3787.2855

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 9442259
$ws.Cells.Item(132, 9).Value = 13517107
$ws.Cells.Item(132, 11).Value = 40551321
$ws.Cells.Item(132, 13).Value = -40548791
$ws.Cells.Item(136, 8).Value = 17565232
$ws.Cells.Item(136, 9).Value = 26316928
$ws.Cells.Item(136, 10).Value = 61840.156
$ws.Cells.Item(136, 11).Value = 78950784
$ws.Cells.Item(136, 12).Value = 185520.468
$ws.Cells.Item(136, 13).Value = -78948234
$ws.Cells.Item(136, 14).Value = -190620.468
